# Update NATMI LR-pair (Nppa-Npr3) data with new TPM-derived values.
# The sending/target clusters now include "ECs" in addition to "MuSCs",
# so the 3-row table (MuSCs x {ECs,FAPs,MuSCs}) becomes a 6-row table
# (ECs x {ECs,FAPs,MuSCs}) followed by (MuSCs x {ECs,FAPs,MuSCs}), with
# refreshed numeric columns throughout.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Nppa"
$ws.Cells.Item(2, 3).Value = "Npr3"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 2
$ws.Cells.Item(2, 6).Value = 0.6666666666666666
$ws.Cells.Item(2, 7).Value = 0.1893306666666667
$ws.Cells.Item(2, 8).Value = 0.567992
$ws.Cells.Item(2, 9).Value = 0.630816252187897
$ws.Cells.Item(2, 10).Value = 0.630816252187897
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 0.912887
$ws.Cells.Item(2, 14).Value = 2.738661
$ws.Cells.Item(2, 15).Value = 0.3341311666818157
$ws.Cells.Item(2, 16).Value = 0.3341311666818156
$ws.Cells.Item(2, 17).Value = 0.1728375043013334
$ws.Cells.Item(2, 18).Value = 1.555537538712
$ws.Cells.Item(2, 19).Value = 0.2107753703053925
$ws.Cells.Item(2, 20).Value = 0.2107753703053925

# Row 3
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Nppa"
$ws.Cells.Item(3, 3).Value = "Npr3"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 2
$ws.Cells.Item(3, 6).Value = 0.6666666666666666
$ws.Cells.Item(3, 7).Value = 0.1893306666666667
$ws.Cells.Item(3, 8).Value = 0.567992
$ws.Cells.Item(3, 9).Value = 0.630816252187897
$ws.Cells.Item(3, 10).Value = 0.630816252187897
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 1.302136
$ws.Cells.Item(3, 14).Value = 3.906408
$ws.Cells.Item(3, 15).Value = 0.4766024939104103
$ws.Cells.Item(3, 16).Value = 0.4766024939104103
$ws.Cells.Item(3, 17).Value = 0.2465342769706667
$ws.Cells.Item(3, 18).Value = 2.218808492736
$ws.Cells.Item(3, 19).Value = 0.3006485989919701
$ws.Cells.Item(3, 20).Value = 0.3006485989919701

# Row 4
$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Nppa"
$ws.Cells.Item(4, 3).Value = "Npr3"
$ws.Cells.Item(4, 4).Value = "MuSCs"
$ws.Cells.Item(4, 5).Value = 2
$ws.Cells.Item(4, 6).Value = 0.6666666666666666
$ws.Cells.Item(4, 7).Value = 0.1893306666666667
$ws.Cells.Item(4, 8).Value = 0.567992
$ws.Cells.Item(4, 9).Value = 0.630816252187897
$ws.Cells.Item(4, 10).Value = 0.630816252187897
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 0.5170986666666667
$ws.Cells.Item(4, 14).Value = 1.551296
$ws.Cells.Item(4, 15).Value = 0.1892663394077741
$ws.Cells.Item(4, 16).Value = 0.189266339407774
$ws.Cells.Item(4, 17).Value = 0.09790263529244446
$ws.Cells.Item(4, 18).Value = 0.881123717632
$ws.Cells.Item(4, 19).Value = 0.1193922828905345
$ws.Cells.Item(4, 20).Value = 0.1193922828905345

# Row 5
$ws.Cells.Item(5, 1).Value = "MuSCs"
$ws.Cells.Item(5, 2).Value = "Nppa"
$ws.Cells.Item(5, 3).Value = "Npr3"
$ws.Cells.Item(5, 4).Value = "ECs"
$ws.Cells.Item(5, 5).Value = 2
$ws.Cells.Item(5, 6).Value = 0.6666666666666666
$ws.Cells.Item(5, 7).Value = 0.1108053333333333
$ws.Cells.Item(5, 8).Value = 0.332416
$ws.Cells.Item(5, 9).Value = 0.3691837478121029
$ws.Cells.Item(5, 10).Value = 0.3691837478121029
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 0.912887
$ws.Cells.Item(5, 14).Value = 2.738661
$ws.Cells.Item(5, 15).Value = 0.3341311666818157
$ws.Cells.Item(5, 16).Value = 0.3341311666818156
$ws.Cells.Item(5, 17).Value = 0.1011527483306667
$ws.Cells.Item(5, 18).Value = 0.910374734976
$ws.Cells.Item(5, 19).Value = 0.1233557963764232
$ws.Cells.Item(5, 20).Value = 0.1233557963764231

# Row 6
$ws.Cells.Item(6, 1).Value = "MuSCs"
$ws.Cells.Item(6, 2).Value = "Nppa"
$ws.Cells.Item(6, 3).Value = "Npr3"
$ws.Cells.Item(6, 4).Value = "FAPs"
$ws.Cells.Item(6, 5).Value = 2
$ws.Cells.Item(6, 6).Value = 0.6666666666666666
$ws.Cells.Item(6, 7).Value = 0.1108053333333333
$ws.Cells.Item(6, 8).Value = 0.332416
$ws.Cells.Item(6, 9).Value = 0.3691837478121029
$ws.Cells.Item(6, 10).Value = 0.3691837478121029
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 1.302136
$ws.Cells.Item(6, 14).Value = 3.906408
$ws.Cells.Item(6, 15).Value = 0.4766024939104103
$ws.Cells.Item(6, 16).Value = 0.4766024939104103
$ws.Cells.Item(6, 17).Value = 0.1442836135253333
$ws.Cells.Item(6, 18).Value = 1.298552521728
$ws.Cells.Item(6, 19).Value = 0.1759538949184402
$ws.Cells.Item(6, 20).Value = 0.1759538949184402

# Row 7
$ws.Cells.Item(7, 1).Value = "MuSCs"
$ws.Cells.Item(7, 2).Value = "Nppa"
$ws.Cells.Item(7, 3).Value = "Npr3"
$ws.Cells.Item(7, 4).Value = "MuSCs"
$ws.Cells.Item(7, 5).Value = 2
$ws.Cells.Item(7, 6).Value = 0.6666666666666666
$ws.Cells.Item(7, 7).Value = 0.1108053333333333
$ws.Cells.Item(7, 8).Value = 0.332416
$ws.Cells.Item(7, 9).Value = 0.3691837478121029
$ws.Cells.Item(7, 10).Value = 0.3691837478121029
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 0.5170986666666667
$ws.Cells.Item(7, 14).Value = 1.551296
$ws.Cells.Item(7, 15).Value = 0.1892663394077741
$ws.Cells.Item(7, 16).Value = 0.189266339407774
$ws.Cells.Item(7, 17).Value = 0.05729729012622222
$ws.Cells.Item(7, 18).Value = 0.515675611136
$ws.Cells.Item(7, 19).Value = 0.06987405651723953
$ws.Cells.Item(7, 20).Value = 0.06987405651723952
